# Leave card update (2/7/2024 8:56am)
# - Remove the now-unused trailing blank table row (row 64), which
#   shrinks Table1 from A8:K64 to A8:K63 and the sheet dimension from
#   A2:K64 to A2:K63.
# - Refresh the view state (split position / selection).
# - Update the footer: widen the "PREPARED BY" signature line and
#   replace the "CERTIFIED CORRECT BY" name/title.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Delete the last (blank) table row -------------------------------
# Row 64 was an empty trailing row of Table1. Deleting the whole sheet
# row removes it from the grid, shrinks the table (ref A8:K64 -> A8:K63)
# and updates the sheet dimension (A2:K64 -> A2:K63) automatically.
$ws.Rows.Item(64).Delete()

# --- 2. View / selection state -------------------------------------------
$excel.ActiveWindow.SplitRow = 12
$ws.Range("M66").Select()

# --- 3. Footer text updates ------------------------------------------------
$ps = $ws.PageSetup
$ps.LeftFooter   = "`n`nPREPARED BY: _____________________`nDATE: &D, &T"
$ps.CenterFooter = "`n`nCERTIFIED CORRECT BY: &UNANETTE B. SUSA&U`n                                             OIC- HRMO"
$ps.RightFooter  = "Page &P of &N"

Write-Host "Leave card updated."
